# Weekly update for "Fruta, Terminal La Palmera de La Serena - Frutilla"
# Inserts 3 new rows of data (728:730) at the top of the date-ordered block
# that starts at row 728, pushing the previously-existing rows 728:782 down
# to 731:785.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 728 downward by 3 (Excel copies formatting from the row above,
# which is exactly what we want since every data row shares the same layout).
$ws.Rows("728:730").Insert()

# --- New row 728 -----------------------------------------------------
$ws.Cells.Item(728, 1).Value = 8
$ws.Cells.Item(728, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(728, 3).Value = "Coquimbo"
$ws.Cells.Item(728, 4).Value = 44746
$ws.Cells.Item(728, 5).Value = 4
$ws.Cells.Item(728, 6).Value = "Fruta"
$ws.Cells.Item(728, 7).Value = 100101
$ws.Cells.Item(728, 8).Value = "Berries"
$ws.Cells.Item(728, 9).Value = 100112025
$ws.Cells.Item(728, 10).Value = "Frutilla"
$ws.Cells.Item(728, 11).Value = "Sin especificar"
$ws.Cells.Item(728, 12).Value = "Especial"
$ws.Cells.Item(728, 13).Value = 300
$ws.Cells.Item(728, 14).Value = 25000
$ws.Cells.Item(728, 15).Value = 26000
$ws.Cells.Item(728, 16).Value = 25500
$ws.Cells.Item(728, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(728, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(728, 19).Value = 3643
$ws.Cells.Item(728, 20).Value = 7

# --- New row 729 -----------------------------------------------------
$ws.Cells.Item(729, 1).Value = 8
$ws.Cells.Item(729, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(729, 3).Value = "Coquimbo"
$ws.Cells.Item(729, 4).Value = 44746
$ws.Cells.Item(729, 5).Value = 4
$ws.Cells.Item(729, 6).Value = "Fruta"
$ws.Cells.Item(729, 7).Value = 100101
$ws.Cells.Item(729, 8).Value = "Berries"
$ws.Cells.Item(729, 9).Value = 100112025
$ws.Cells.Item(729, 10).Value = "Frutilla"
$ws.Cells.Item(729, 11).Value = "Sin especificar"
$ws.Cells.Item(729, 12).Value = "Primera"
$ws.Cells.Item(729, 13).Value = 400
$ws.Cells.Item(729, 14).Value = 20000
$ws.Cells.Item(729, 15).Value = 21000
$ws.Cells.Item(729, 16).Value = 20500
$ws.Cells.Item(729, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(729, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(729, 19).Value = 2929
$ws.Cells.Item(729, 20).Value = 7

# --- New row 730 -----------------------------------------------------
$ws.Cells.Item(730, 1).Value = 8
$ws.Cells.Item(730, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(730, 3).Value = "Coquimbo"
$ws.Cells.Item(730, 4).Value = 44746
$ws.Cells.Item(730, 5).Value = 4
$ws.Cells.Item(730, 6).Value = "Fruta"
$ws.Cells.Item(730, 7).Value = 100101
$ws.Cells.Item(730, 8).Value = "Berries"
$ws.Cells.Item(730, 9).Value = 100112025
$ws.Cells.Item(730, 10).Value = "Frutilla"
$ws.Cells.Item(730, 11).Value = "Sin especificar"
$ws.Cells.Item(730, 12).Value = "Segunda"
$ws.Cells.Item(730, 13).Value = 240
$ws.Cells.Item(730, 14).Value = 16000
$ws.Cells.Item(730, 15).Value = 17000
$ws.Cells.Item(730, 16).Value = 16500
$ws.Cells.Item(730, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(730, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(730, 19).Value = 2357
$ws.Cells.Item(730, 20).Value = 7
